# Applies the "Weekly" -> "Monthly" Scottish payroll input-sheet rename.
#
# Semantic changes described by the diff:
#   1. Worksheet "GeneralTaxRateWeekly"        -> "GeneralTaxRateMonthly"
#   2. Worksheet "ProcessPayrollForWeeklyTax"  -> "ProcessPayrollForMonthlyTax"
#   3. On sheet "first", the rows that reference those two sheet/testcase
#      names are updated to match.
#   4. The "DO NOT TOUCH AUTOMATION EMP 107" marker text (present on the
#      GeneralTaxRate*, ProcessPayrollFor*Tax and TestReports sheets) is
#      updated to "DO NOT TOUCH AUTOMATION EMP 105".
#   5. The workbook is left with the (renamed) GeneralTaxRateMonthly sheet
#      as the active / selected sheet, and the saved selections on every
#      sheet move to match the diff's view state.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: rename the two sheets ------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("GeneralTaxRateWeekly")
$wsGeneral.Name = "GeneralTaxRateMonthly"

$wsProcess = $wb.Worksheets.Item("ProcessPayrollForWeeklyTax")
$wsProcess.Name = "ProcessPayrollForMonthlyTax"

# --- 3: update the "first" index sheet's TC names to match ----------------------
$wsFirst = $wb.Worksheets.Item("first")
$wsFirst.Range("A3").Value = "GeneralTaxRateMonthly"
$wsFirst.Range("A4").Value = "ProcessPayrollForMonthlyTax"

# --- 4: update the "DO NOT TOUCH AUTOMATION EMP 107" marker everywhere ----------
$wsGeneral.Range("A2").Value = "DO NOT TOUCH AUTOMATION EMP 105"
$wsProcess.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"

$wsReports = $wb.Worksheets.Item("TestReports")
$wsReports.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"

# --- 5: update saved selections / active sheet to match the new view state ------
$wsFirst.Select()
$wsFirst.Range("F5").Select()

$wsProcess.Select()
$wsProcess.Range("C7").Select()

$wsReports.Select()
$wsReports.Range("F7").Select()

$wsGeneral.Select()
$wsGeneral.Range("G11").Select()
